$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / header text updates (shared strings 6 and 9) ---
$ws.Range("A8").Value = "Volume 30   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/12/2023  Through  6/18/2023"

# --- Data table updates (rows 15-29) ---
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Value = "'0"
$ws.Range("E15").NumberFormat = "General"
$ws.Range("E15").Value = "'***.*"
$ws.Range("G15").Value = 1
$ws.Range("N15").Value = -92.682926829268
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = -48
$ws.Range("I16").Value = 104
$ws.Range("J16").Value = 113
$ws.Range("K16").Value = -7.964601769911
$ws.Range("L16").Value = 18.181818181818
$ws.Range("M16").Value = -20
$ws.Range("N16").Value = -85.057471264367
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 41
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 115.789473684211
$ws.Range("I17").Value = 157
$ws.Range("J17").Value = 133
$ws.Range("K17").Value = 18.045112781954
$ws.Range("L17").Value = 26.612903225806
$ws.Range("M17").Value = 72.527472527472
$ws.Range("N17").Value = -62.619047619047
$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 2
$ws.Range("E18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -54.545454545454
$ws.Range("J18").Value = 72
$ws.Range("K18").Value = -5.555555555555
$ws.Range("L18").Value = 33.333333333333
$ws.Range("M18").Value = -15
$ws.Range("N18").Value = -92.953367875647
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -23.076923076923
$ws.Range("F19").Value = 43
$ws.Range("H19").Value = -14
$ws.Range("I19").Value = 236
$ws.Range("J19").Value = 295
$ws.Range("K19").Value = -20
$ws.Range("L19").Value = -9.230769230769
$ws.Range("M19").Value = 61.643835616438
$ws.Range("N19").Value = -56.855575868372
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 23
$ws.Range("H20").Value = -30.434782608695
$ws.Range("I20").Value = 118
$ws.Range("J20").Value = 166
$ws.Range("K20").Value = -28.915662650602
$ws.Range("L20").Value = 40.476190476190
$ws.Range("M20").Value = 210.526315789474
$ws.Range("N20").Value = -88.104838709677
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = 17.391304347826
$ws.Range("F21").Value = 119
$ws.Range("G21").Value = 129
$ws.Range("H21").Value = -7.751937984496
$ws.Range("I21").Value = 687
$ws.Range("J21").Value = 792
$ws.Range("K21").Value = -13.257575757575
$ws.Range("L21").Value = 10.628019323671
$ws.Range("M21").Value = 38.229376257545
$ws.Range("N21").Value = -81.402273957769
$ws.Range("C22").NumberFormat = "General"
$ws.Range("C22").Value = "'0"
$ws.Range("M22").Value = 76.923076923076
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Value = "'0"
$ws.Range("E23").NumberFormat = "General"
$ws.Range("E23").Value = "'***.*"
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = -80
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 133.333333333333
$ws.Range("F24").Value = 170
$ws.Range("G24").Value = 96
$ws.Range("H24").Value = 77.083333333333
$ws.Range("I24").Value = 650
$ws.Range("J24").Value = 587
$ws.Range("K24").Value = 10.732538330494
$ws.Range("L24").Value = 71.052631578947
$ws.Range("M24").Value = 130.496453900709
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 128.571428571429
$ws.Range("I25").Value = 259
$ws.Range("J25").Value = 207
$ws.Range("K25").Value = 25.120772946859
$ws.Range("L25").Value = 48
$ws.Range("M25").Value = 5.284552845528
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Value = "'0"
$ws.Range("E26").NumberFormat = "General"
$ws.Range("E26").Value = "'***.*"
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 0
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 1
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Value = "'0"
$ws.Range("E27").NumberFormat = "General"
$ws.Range("E27").Value = "'***.*"
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 24
$ws.Range("K27").Value = 4.347826086956
$ws.Range("L27").Value = -7.692307692307
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Value = "'0"
$ws.Range("E28").NumberFormat = "General"
$ws.Range("E28").Value = "'***.*"
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = -50
$ws.Range("N28").Value = -97.247706422018
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Value = "'0"
$ws.Range("E29").NumberFormat = "General"
$ws.Range("E29").Value = "'***.*"
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -50
$ws.Range("N29").Value = -96.907216494845
